$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add E4 as an empty cell sharing the same "key" formatting as C4/D4.
# ---------------------------------------------------------------------------
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Clear the fill from E13:I14 (was the green "key" fill, now no fill).
#    Doing this first creates the "fill cleared" cell style before the new
#    font-colour style below, matching the original authoring order.
# ---------------------------------------------------------------------------
$ws.Range("E13:I14").Interior.Pattern = -4142

# ---------------------------------------------------------------------------
# 3. New row 20 - "log(arcsin(x))" transformation, highlighted in the
#    accent1 theme colour font.
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "log(arcsin(x))"
$ws.Range("B20").Value = 0.97626000000000002
$ws.Range("C20").Value = "<2.2 e^-16"
$ws.Range("D20").Value = "slight right-skew"
$ws.Range("E20").Value = "Not great, deviates at bottom big curve at top"
$ws.Range("A20:E20").Font.ThemeColor = 5

# ---------------------------------------------------------------------------
# 4. New row 21 - "ln(sqrt(x))" transformation (proportional-data columns).
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "ln(sqrt(x))"
$ws.Range("F21").Value = 0.97814999999999996
$ws.Range("G21").Value = "1.068 e^-12"
$ws.Range("H21").Value = "Slight left-skew and lumpiness"
$ws.Range("I21").Value = "A bit skewiff at the lower values, but pretty nice"
$ws.Range("E13").Copy()
$ws.Range("F21:I21").PasteSpecial(-4122)
$ws.Range("F21").Value = 0.97814999999999996
$ws.Range("G21").Value = "1.068 e^-12"
$ws.Range("H21").Value = "Slight left-skew and lumpiness"
$ws.Range("I21").Value = "A bit skewiff at the lower values, but pretty nice"

# ---------------------------------------------------------------------------
# 5. New row 22 - "cuberoot(arcsin)" transformation, highlighted in the same
#    blue font used elsewhere for favourite candidates.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("A22:E22").PasteSpecial(-4122)
$ws.Range("A22").Value = "cuberoot(arcsin)"
$ws.Range("B22").Value = 0.92232999999999998
$ws.Range("C22").Value = "<2.2 e^-16"
$ws.Range("D22").Value = "left skew"
$ws.Range("E22").Value = "S-shaped (less so at top)"

# ---------------------------------------------------------------------------
# 6. New row 23 - "ln(x^2)" transformation (proportional-data columns).
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = "ln(x^2)"
$ws.Range("E13").Copy()
$ws.Range("F23:I23").PasteSpecial(-4122)
$ws.Range("F23").Value = 0.97814999999999996
$ws.Range("G23").Value = "1.068 e^-12"
$ws.Range("H23").Value = "Slight left-skew and lumpiness"
$ws.Range("I23").Value = "A bit skewiff at the lower values, but pretty nice"

# ---------------------------------------------------------------------------
# 7. New row 24 - "sqrt(ln())" transformation, keeps the original green
#    "key" fill formatting (same style as before row 13/14 lost theirs).
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = "sqrt(ln())"
$ws.Range("D4").Copy()
$ws.Range("F24:I24").PasteSpecial(-4122)
$ws.Range("F24").Value = 0.98341000000000001
$ws.Range("G24").Value = "6.249 e^-10"
$ws.Range("H24").Value = "Pretty good, lump on left"
$ws.Range("I24").Value = "Pretty straight, but with a tail at the bottom"

# ---------------------------------------------------------------------------
# 8. New row 25 - "cbrt(ln())" transformation (proportional-data columns).
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "cbrt(ln())"
$ws.Range("E13").Copy()
$ws.Range("F25:I25").PasteSpecial(-4122)
$ws.Range("F25").Value = 0.97814999999999996
$ws.Range("G25").Value = "1.68 e^-12"
$ws.Range("H25").Value = "Slight left-skew and lumpiness"
$ws.Range("I25").Value = "Pretty straight"

# ---------------------------------------------------------------------------
# 9. Park the selection where the author left it.
# ---------------------------------------------------------------------------
$ws.Range("E29").Select()
